# Fruta / hortaliza, semanal
# Insert a new weekly record at row 320 (pushing the existing rows 320-349
# down to 321-350) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 320, shifting everything
# from row 320 downward (previously ending at row 349) down by one row,
# so the sheet now spans to row 350.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with this week's data.
$ws.Cells.Item(320, 1).Value()  = 10
$ws.Cells.Item(320, 2).Value()  = "Vega Modelo de Temuco"
$ws.Cells.Item(320, 3).Value()  = "La Araucanía"
$ws.Cells.Item(320, 4).Value()  = 44826
$ws.Cells.Item(320, 5).Value()  = 9
$ws.Cells.Item(320, 6).Value()  = 100112001
$ws.Cells.Item(320, 7).Value()  = "Berenjena"
$ws.Cells.Item(320, 8).Value()  = "Sin especificar"
$ws.Cells.Item(320, 9).Value()  = "Primera"
$ws.Cells.Item(320, 10).Value() = 100
$ws.Cells.Item(320, 11).Value() = 15000
$ws.Cells.Item(320, 12).Value() = 15000
$ws.Cells.Item(320, 13).Value() = 15000
$ws.Cells.Item(320, 14).Value() = "$/caja 40 unidades"
$ws.Cells.Item(320, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(320, 16).Value() = 375
$ws.Cells.Item(320, 17).Value() = 40
$ws.Cells.Item(320, 18).Value() = "Hortaliza"
